# Applies the "minor changes in manual csv file" edit to Sheet1.
#
# Summary of the edit (from the OOXML diff):
#  - New column J header "FAIR-Test-All" (J1) plus a few "No"/"no" flags
#    further down column J (rows 29, 34, 38, 40, 41).
#  - Row 8's repository-link cell (B8) text was replaced.
#  - New column M header "Included or Not" (M1) with "in" filled down
#    column M for most data rows.
#  - A new cell P27 containing the literal text "=".
#  - B24 and B28 (which already held plain-text URLs) became real
#    hyperlinks, picking up the workbook's "Hyperlink" cell style.
#  - The active selection moved to P27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column J header ------------------------------------------------
$ws.Range("J1").Value = "FAIR-Test-All"

# --- Fix the broken-looking repository link text in B8 ------------------
$ws.Range("B8").Value = "https://github.com/SemanticM+B9+B10+B11+B+B23"

# --- New column M header + "in" filled down the data rows ---------------
$ws.Range("M1").Value = "Included or Not"

$inRows = @(2,3,4,5,6,7,8,9,10,11,12,13,29,34,38,40,41)
foreach ($r in $inRows) {
    $ws.Cells.Item($r, 13).Value = "in"   # column M = 13
}

# --- "No" / "no" flags in column J ---------------------------------------
foreach ($r in @(29,34)) {
    $ws.Cells.Item($r, 10).Value = "No"   # column J = 10
}
foreach ($r in @(38,40,41)) {
    $ws.Cells.Item($r, 10).Value = "no"
}

# --- P27 literal "=" text (force text so it isn't parsed as a formula) --
$p27 = $ws.Range("P27")
$p27.Value = "'="
$p27.Style = "Normal"

# --- Turn the existing plain-text URLs in B24 / B28 into real hyperlinks
$ws.Hyperlinks.Add($ws.Range("B24"), "https://github.com/fastscape-lem") | Out-Null
$ws.Range("B24").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("B28"), "https://github.com/GFK-in-Potsdam") | Out-Null
$ws.Range("B28").Style = "Hyperlink"

# --- Match the saved selection/active cell -------------------------------
$ws.Range("P27").Select() | Out-Null
